$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Ref,
        [string]$NewValue,
        [bool]$ForceText = $false
    )
    $cell = $ws.Range($Ref)
    if ($ForceText) {
        # Leading apostrophe forces Excel to keep a numeric-looking string as
        # literal text instead of silently coercing it to a Double. This
        # stamps a quotePrefix style on the cell, so clear formatting right
        # back to the (unstyled) default afterwards to avoid leaving any
        # style residue behind.
        $cell.Value = "'" + $NewValue
        $cell.ClearFormats()
    } else {
        $cell.Value = $NewValue
    }
}

# Row 2 - Bitcoin
Set-TextValue "D2" "29.061.44" $false
Set-TextValue "E2" "  +0.92%  " $false

# Row 3 - Ethereum
Set-TextValue "D3" "1.898.05" $false
Set-TextValue "E3" "  +0.89%  " $false

# Row 4 - TetherUSD
Set-TextValue "D4" "1.004" $true
Set-TextValue "E4" "  -0.29%  " $false

# Row 5 - BNB
Set-TextValue "D5" "327.28" $true
Set-TextValue "E5" "  +0.91%  " $false

# Row 6 - USDC
Set-TextValue "E6" "  -0.21%  " $false

# Row 7 - XRP
Set-TextValue "D7" "0.4648" $true
Set-TextValue "E7" "  -0.64%  " $false

# Row 8 - Cardano
Set-TextValue "D8" "0.3918" $true
Set-TextValue "E8" "  -0.42%  " $false

# Row 9 - OKB
Set-TextValue "D9" "47.03" $true
Set-TextValue "E9" "  +1.02%  " $false

# Row 10 - Dogecoin
Set-TextValue "D10" "0.07950" $true
Set-TextValue "E10" "  +0.20%  " $false

# Row 11 - Polygon
Set-TextValue "D11" "1.012" $true
Set-TextValue "E11" "  +3.45%  " $false

# Row 12 - Solana
Set-TextValue "D12" "22.12" $true
Set-TextValue "E12" "  -0.97%  " $false

# Row 13 - WrappedEther
Set-TextValue "D13" "1.921.30" $false
Set-TextValue "E13" "  +0.17%  " $false

# Row 14 - Chainlink
Set-TextValue "D14" "7.130" $true
Set-TextValue "E14" "  +1.57%  " $false

# Row 15 - Polkadot
Set-TextValue "D15" "5.779" $true
Set-TextValue "E15" "  +0.67%  " $false

# Row 16 - TRON
Set-TextValue "D16" "0.06976" $true

# Row 17 - Litecoin
Set-TextValue "D17" "89.20" $true
Set-TextValue "E17" "  +0.56%  " $false

# Row 18 - BinanceUSD
Set-TextValue "D18" "1.003" $true
Set-TextValue "E18" "  -0.28%  " $false

# Row 19 - ShibaInu
Set-TextValue "D19" "0.00001014" $true
Set-TextValue "E19" "  +0.40%  " $false

# Row 20 - Avalanche
Set-TextValue "D20" "17.28" $true
Set-TextValue "E20" "  +1.79%  " $false

# Row 21 - Dai
Set-TextValue "D21" "1.004" $true
Set-TextValue "E21" "  -0.06%  " $false

# Row 22 - WrappedBTC
Set-TextValue "D22" "29.070.54" $false
Set-TextValue "E22" "  +0.88%  " $false

# Row 23 - Uniswap
Set-TextValue "D23" "5.358" $true
Set-TextValue "E23" "  +0.22%  " $false

# Row 24 - Cosmos
Set-TextValue "D24" "11.12" $true
Set-TextValue "E24" "  +0.29%  " $false

# Row 25 - WrappedliquidstakedEther2.0
Set-TextValue "D25" "2.141.86" $false
Set-TextValue "E25" "  +0.11%  " $false

# Row 26 - Toncoin
Set-TextValue "E26" "  -2.99%  " $false

# Row 27 - Monero
Set-TextValue "D27" "155.26" $true
Set-TextValue "E27" "  +0.98%  " $false

# Row 28 - EthereumClassic
Set-TextValue "D28" "19.86" $true
Set-TextValue "E28" "  +2.35%  " $false

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "5.881" $true
Set-TextValue "E29" "  +2.19%  " $false

# Row 30 - LidoDAOToken
Set-TextValue "E30" "  -0.43%  " $false

# Row 31 - BitcoinCash
Set-TextValue "D31" "119.74" $true
Set-TextValue "E31" "  -0.26%  " $false

# Row 32 - Stellar
Set-TextValue "D32" "0.09391" $true
Set-TextValue "E32" "  -0.06%  " $false

# Row 33 - ImmutableX
Set-TextValue "E33" "  -0.08%  " $false

# Row 34 - Filecoin
Set-TextValue "D34" "5.367" $true
Set-TextValue "E34" "  +0.97%  " $false

# Row 35 - ARBITRUM
Set-TextValue "D35" "1.352" $true
Set-TextValue "E35" "  -0.31%  " $false

# Row 36 - HuobiToken
Set-TextValue "E36" "  -2.88%  " $false

# Row 37 - Hedera
Set-TextValue "D37" "0.05837" $true
Set-TextValue "E37" "  -1.36%  " $false

# Row 38 - TrustWalletToken
Set-TextValue "E38" "  +1.28%  " $false

# Row 39 - FraxShare
Set-TextValue "D39" "8.103" $true
Set-TextValue "E39" "  +2.56%  " $false

# Row 40 - VeChain
Set-TextValue "D40" "0.02096" $true
Set-TextValue "E40" "  -1.29%  " $false

# Row 41 - TheSandbox
Set-TextValue "D41" "0.5836" $true
Set-TextValue "E41" "  +2.00%  " $false

# Row 43 - Algorand
Set-TextValue "D43" "0.1816" $true
Set-TextValue "E43" "  +1.04%  " $false

# Row 44 - Aptos
Set-TextValue "D44" "9.999" $true
Set-TextValue "E44" "  +0.00%  " $false

# Row 45 - was Decentraland, now RenderToken
Set-TextValue "B45" "RenderToken" $false
Set-TextValue "C45" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" $false
Set-TextValue "D45" "2.262" $true
Set-TextValue "E45" "  +7.38%  " $false

# Row 46 - was RenderToken, now Decentraland
Set-TextValue "B46" "Decentraland" $false
Set-TextValue "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" $false
Set-TextValue "D46" "0.5465" $true
Set-TextValue "E46" "  +2.21%  " $false

# Row 47 - EnergySwap
Set-TextValue "D47" "11.98" $true
Set-TextValue "E47" "  +1.30%  " $false

# Row 48 - Cronos
Set-TextValue "D48" "0.07205" $true
Set-TextValue "E48" "  -1.74%  " $false

# Row 49 - NEARProtocol
Set-TextValue "D49" "1.870" $true
Set-TextValue "E49" "  +1.21%  " $false

# Row 50 - WEMIXToken
Set-TextValue "D50" "1.127" $true
Set-TextValue "E50" "  -2.13%  " $false

# Row 51 - Quant
Set-TextValue "D51" "113.41" $true
Set-TextValue "E51" "  -0.66%  " $false
